$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'301.68"
$ws.Range("E2").Value = "'-2.03%"
$ws.Range("D3").Value = "'40.56"
$ws.Range("E3").Value = "'-1.23%"
$ws.Range("D4").Value = "'5.147"
$ws.Range("E4").Value = "'-1.32%"
$ws.Range("D5").Value = "'0.07574"
$ws.Range("E5").Value = "'-1.21%"
$ws.Range("D6").Value = "'4.338"
$ws.Range("E6").Value = "'0.62%"
$ws.Range("D7").Value = "'1.652"
$ws.Range("E7").Value = "'0.73%"
$ws.Range("D8").Value = "'0.9420"
$ws.Range("E8").Value = "'2.98%"
$ws.Range("D9").Value = "'0.1200"
$ws.Range("E9").Value = "'-0.89%"
$ws.Range("D10").Value = "'0.1793"
$ws.Range("E10").Value = "'-1.42%"
$ws.Range("D11").Value = "'0.08976"
$ws.Range("E11").Value = "'-2.20%"
$ws.Range("D12").Value = "'0.04167"
$ws.Range("E12").Value = "'-0.43%"
$ws.Range("E13").Value = "'0.30%"
$ws.Range("D14").Value = "'0.001291"
$ws.Range("E14").Value = "'2.57%"
$ws.Range("E15").Value = "'-0.24%"
$ws.Range("D16").Value = "'3.342"
$ws.Range("E16").Value = "'0.05%"
$ws.Range("D17").Value = "'2.424"
$ws.Range("E17").Value = "'-0.21%"
$ws.Range("D19").Value = "'7.634"
$ws.Range("E19").Value = "'3.14%"
$ws.Range("D20").Value = "'0.1355"
$ws.Range("E20").Value = "'-1.99%"
$ws.Range("D21").Value = "'0.2810"
$ws.Range("E21").Value = "'3.52%"
$ws.Range("D22").Value = "'0.03866"
$ws.Range("E22").Value = "'-3.26%"
$ws.Range("D23").Value = "'0.001285"
$ws.Range("E23").Value = "'2.11%"
$ws.Range("D24").Value = "'0.003972"
$ws.Range("E24").Value = "'-9.31%"
$ws.Range("D25").Value = "'0.0001303"
$ws.Range("E25").Value = "'0.09%"
$ws.Range("D26").Value = "'0.0003735"
$ws.Range("E26").Value = "'-95.03%"
$ws.Range("D38").Value = "'0.02360"
$ws.Range("E38").Value = "'-4.89%"
$ws.Range("D39").Value = "'0.05114"
$ws.Range("E39").Value = "'-4.01%"
$ws.Range("D40").Value = "'0.007710"
$ws.Range("E40").Value = "'-1.75%"
$ws.Range("E41").Value = "'-1.21%"
$ws.Range("D42").Value = "'0.007589"
$ws.Range("E42").Value = "'16.58%"
$ws.Range("D43").Value = "'0.003686"
$ws.Range("E43").Value = "'92.74%"
$ws.Range("D44").Value = "'0.007403"
$ws.Range("E44").Value = "'-10.37%"
$ws.Range("D45").Value = "'0.3263"
$ws.Range("E45").Value = "'-2.49%"
$ws.Range("D46").Value = "'0.00006821"
$ws.Range("E46").Value = "'1.63%"
$ws.Range("D47").Value = "'0.00000000753"
$ws.Range("E47").Value = "'0.21%"
$ws.Range("D48").Value = "'0.2604"
$ws.Range("E48").Value = "'-33.43%"
$ws.Range("D49").Value = "'0.004214"
$ws.Range("E49").Value = "'35.75%"
$ws.Range("D50").Value = "'0.00002107"
$ws.Range("E50").Value = "'0.21%"
$ws.Range("D51").Value = "'0.0002007"
$ws.Range("E51").Value = "'0.21%"
